$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.497.23"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.457.17"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'583.95"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").Value = "'143.91"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "2.453.99"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "'5.23"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'26.56"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "2.897.96"
$ws.Range("D17").Value = "62.208.62"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "2.453.41"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "'10.88"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "'329.36"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'1.97"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'65.83"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "'9.26"
$ws.Range("E26").Value = "  +4.60%  "
$ws.Range("D27").Value = "'592.59"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").Value = "0.0₃0981"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("D29").Value = "2.574.78"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "'0.136"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'0.380"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").Value = "'154.49"
$ws.Range("E39").Value = "  +5.57%  "
$ws.Range("D40").Value = "'5.34"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "'18.48"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'43.04"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "'143.63"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "0.0₆0259"
$ws.Range("E48").Value = "  +17.90%  "
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "'0.0526"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "'20.01"
$ws.Range("E51").Value = "  -0.70%  "
